$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update effort/sprint for "12. Prikaz obavijesti za uzimanje terapije" (row 16) ---
$ws.Range("D16").Value = 16
$ws.Range("F16").Value = 3

# --- Insert a new backlog item before the last row ------------------------
# The old last row (15. Bar-code čitač, row 19) moves down to row 20, and a
# brand new row 19 (15. Izmjena korisničkih podataka) takes its place; the
# old row then gets renumbered to "16.".
#
# Move row 19's current content+formatting down to row 20 first.
$ws.Range("B19:F19").Copy($ws.Range("B20:F20"))

# Give the (soon to be new) row 19 the same formatting as row 18 (a normal,
# non-final data row), matching Excel's default "insert row" behaviour.
$ws.Range("B18:F18").Copy($ws.Range("B19:F19"))

# Helper cell (well outside the used range) used only to coerce literal,
# dot-terminated numbering labels ("15.", "16.") to be stored as text
# instead of being auto-converted to numbers - mirrors how every other
# "N." label in this sheet is stored as a shared string.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$helper.Value = "15."
$helper.Copy()
$ws.Range("B19").PasteSpecial(-4163)  # xlPasteValues - keep row19's own style

$helper.Value = "16."
$helper.Copy()
$ws.Range("B20").PasteSpecial(-4163)  # xlPasteValues - keep row20's own style

$helper.Clear()

# Fill in the rest of the new row 19 content.
$ws.Range("C19").Value = "Izmjena korisničkih podataka"
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 9
$ws.Range("F19").Value = 2

# --- Match the target workbook's view/selection state ---------------------
$ws.Range("I17").Select()
